$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I64").Value = -5.9
$ws.Range("J64").Value = 10.8
$ws.Range("K64").Value = 19.2
$ws.Range("M64").Value = 5.5
$ws.Range("N64").Value = 8.9
$ws.Range("Q64").Value = 1900
$ws.Range("U64").Value = -0.26
$ws.Range("V64").Value = 0.25
$ws.Range("W64").Value = 46.5
$ws.Range("X64").Value = -7.8
$ws.Range("Y64").Value = 9.1
$ws.Range("AA64").Value = -3.2
$ws.Range("AP64").Value = -7.5
$ws.Range("AQ64").Value = 8.699999999999999
$ws.Range("H65").Value = 11.7
$ws.Range("I65").Value = -1.7
$ws.Range("K65").Value = 8.4
$ws.Range("L65").Value = -1.7
$ws.Range("M65").Value = 1.3
$ws.Range("N65").Value = 17.4
$ws.Range("O65").Value = -0.3
$ws.Range("Q65").Value = 640
$ws.Range("R65").Value = -260
$ws.Range("S65").Value = 200
$ws.Range("X65").Value = -0.95
$ws.Range("Y65").Value = 1.64
$ws.Range("Z65").Value = 8.56
$ws.Range("AA65").Value = -0.28
$ws.Range("AB65").Value = 0.41
$ws.Range("AO65").Value = 19.18
$ws.Range("AP65").Value = -0.93
$ws.Range("AQ65").Value = 1.71
$ws.Range("I67").Value = -2.2
$ws.Range("X67").Value = -0.93
$ws.Range("Y67").Value = 2.7
$ws.Range("AP67").Value = -0.93
$ws.Range("AQ67").Value = 2.81
$ws.Range("S68").Value = 240
$ws.Range("Z68").Value = 4.31
$ws.Range("L72").Value = -0.28
$ws.Range("V72").Value = 0.23
$ws.Range("Y72").Value = 1.7
$ws.Range("I74").Value = -3.1
$ws.Range("K74").Value = 29
$ws.Range("M74").Value = 3.3
$ws.Range("Q74").Value = 890
$ws.Range("R74").Value = -370
$ws.Range("S74").Value = 260
$ws.Range("W74").Value = 63.5
$ws.Range("X74").Value = -3.4
$ws.Range("Y74").Value = 4.5
$ws.Range("AO74").Value = 60.2
$ws.Range("AP74").Value = -3.2
$ws.Range("AQ74").Value = 4.1
$ws.Range("H76").Value = 37.7
$ws.Range("J76").Value = 9.300000000000001
$ws.Range("L76").Value = -7.3
$ws.Range("M76").Value = 6.3
$ws.Range("S76").Value = 1020
$ws.Range("V76").Value = 0.21
$ws.Range("W76").Value = 65.09999999999999
$ws.Range("Y76").Value = 8.1
$ws.Range("AB76").Value = 3.7
$ws.Range("AQ76").Value = 7.5
$ws.Range("J77").Value = 103
$ws.Range("L77").Value = -5.5
$ws.Range("M77").Value = 9.199999999999999
$ws.Range("R77").Value = -2000
$ws.Range("U77").Value = -0.46
$ws.Range("V77").Value = 0.42
$ws.Range("X77").Value = -26
$ws.Range("Z77").Value = 19.8
$ws.Range("AA77").Value = -5.2
$ws.Range("AB77").Value = 10.5
$ws.Range("AG77").Value = -0.29
$ws.Range("AH77").Value = 0.53
$ws.Range("AP77").Value = -26
$ws.Range("AQ77").Value = 99
$ws.Range("J84").Value = 6.7
$ws.Range("K84").Value = 32.7
$ws.Range("M84").Value = 4.8
$ws.Range("R84").Value = -650
$ws.Range("S84").Value = 500
$ws.Range("U84").Value = -0.16
$ws.Range("W84").Value = 72.7
$ws.Range("X84").Value = -5.3
$ws.Range("AA84").Value = -2.7
$ws.Range("AB84").Value = 3.3
$ws.Range("AO84").Value = 68.7
$ws.Range("AP84").Value = -4.8
$ws.Range("AQ84").Value = 6.7
$ws.Range("H88").Value = 60
$ws.Range("I88").Value = -29
$ws.Range("J88").Value = 166
$ws.Range("K88").Value = 24
$ws.Range("L88").Value = -13
$ws.Range("M88").Value = 36
$ws.Range("N88").Value = 4.7
$ws.Range("O88").Value = -2.9
$ws.Range("P88").Value = 2.5
$ws.Range("Q88").Value = 7100
$ws.Range("R88").Value = -4400
$ws.Range("S88").Value = 13900
$ws.Range("T88").Value = 0.16
$ws.Range("U88").Value = -0.49
$ws.Range("V88").Value = 0.58
$ws.Range("W88").Value = 84
$ws.Range("X88").Value = -48
$ws.Range("Y88").Value = 169
$ws.Range("Z88").Value = 34
$ws.Range("AA88").Value = -18
$ws.Range("AB88").Value = 44
$ws.Range("AF88").Value = 1.04
$ws.Range("AG88").Value = -0.57
$ws.Range("AH88").Value = 1.47
$ws.Range("AO88").Value = 88
$ws.Range("AP88").Value = -47
$ws.Range("AQ88").Value = 169
$ws.Range("M90").Value = 2
$ws.Range("S90").Value = 480
$ws.Range("H91").Value = 38
$ws.Range("I91").Value = -22
$ws.Range("J91").Value = 130
$ws.Range("K91").Value = 11.3
$ws.Range("L91").Value = -6
$ws.Range("M91").Value = 24.3
$ws.Range("N91").Value = 4.5
$ws.Range("O91").Value = -3
$ws.Range("P91").Value = 2.7
$ws.Range("Q91").Value = 3500
$ws.Range("R91").Value = -2200
$ws.Range("S91").Value = 12500
$ws.Range("T91").Value = 0.27
$ws.Range("U91").Value = -0.58
$ws.Range("V91").Value = 0.54
$ws.Range("W91").Value = 49.3
$ws.Range("X91").Value = -22
$ws.Range("Y91").Value = 132
$ws.Range("Z91").Value = 15
$ws.Range("AA91").Value = -4
$ws.Range("AB91").Value = 29.5
$ws.Range("AF91").Value = 0.59
$ws.Range("AG91").Value = -0.32
$ws.Range("AH91").Value = 1.43
$ws.Range("AO91").Value = 48
$ws.Range("AP91").Value = -22
$ws.Range("AQ91").Value = 132
